# "Add trustLevel in flow"
# EntryPoint (sheet1 / Table1): insert a new row for the PROM_PROC /
# Administration entry point, and raise the existing rows' Trust Level
# from "Unknown" to "Operational".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new worksheet row right under the header row; this shifts the
# existing PROC_RABBIT / PROC_KAFKA / MINIO_PROC rows down by one.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "PROM_PROC"
$ws.Range("B2").Value = "PROM_PROC"
$ws.Range("C2").Value = "Publish message"
$ws.Range("D2").Value = "Administration"
$ws.Range("E2").Value = "MyProcess"

# Grow Table1 so it covers the newly inserted row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E5"))

# The previously-existing rows (now rows 3-5) move from "Unknown" to
# "Operational" for their Trust Level.
$ws.Range("D3").Value = "Operational"
$ws.Range("D4").Value = "Operational"
$ws.Range("D5").Value = "Operational"

# Column D widened (bestFit) to accommodate the new "Administration" /
# "Operational" text.
$ws.Columns.Item(4).ColumnWidth = 13.43

# Trust Boundaries sheet (sheet2 / Table2): the Enablers / Enablers2 rows
# swap order.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Enablers2"
$ws2.Range("A3").Value = "Enablers"
